$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 1.381856437894394
$ws.Range("F1").Value = -1.570796384046431

$ws.Range("E2").Value = 1.382896364020461
$ws.Range("F2").Value = -1.570796384423863

$ws.Range("E3").Value = 1.387555922519023
$ws.Range("F3").Value = -1.570796386115009

$ws.Range("E4").Value = 1.394111478926988
$ws.Range("F4").Value = -1.570796388494291

$ws.Range("E5").Value = 1.39877103742555
$ws.Range("F5").Value = -1.570796390185437

$ws.Range("E6").Value = 1.399810963551616
$ws.Range("F6").Value = -1.570796390562869
